$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 11 rows of data (rows 32-42), shifting dimension from A1:B42 to A1:B31
$ws.Rows("32:42").Delete()

# Update remaining data rows (2-31) with the new winter-DST-aligned values
$ws.Range("A2").Value = 5450
$ws.Range("B2").Value = 45954
$ws.Range("A3").Value = 5367
$ws.Range("B3").Value = 45954.01041666666
$ws.Range("A4").Value = 5340
$ws.Range("B4").Value = 45954.02083333334
$ws.Range("A5").Value = 5343
$ws.Range("B5").Value = 45954.03125
$ws.Range("A6").Value = 5286
$ws.Range("B6").Value = 45954.04166666666
$ws.Range("A7").Value = 5269
$ws.Range("B7").Value = 45954.05208333334
$ws.Range("A8").Value = 5280
$ws.Range("B8").Value = 45954.0625
$ws.Range("A9").Value = 5232
$ws.Range("B9").Value = 45954.07291666666
$ws.Range("A10").Value = 5194
$ws.Range("B10").Value = 45954.08333333334
$ws.Range("A11").Value = 5214
$ws.Range("B11").Value = 45954.09375
$ws.Range("A12").Value = 5249
$ws.Range("B12").Value = 45954.10416666666
$ws.Range("A13").Value = 5290
$ws.Range("B13").Value = 45954.11458333334
$ws.Range("A14").Value = 5263
$ws.Range("B14").Value = 45954.125
$ws.Range("A15").Value = 5224
$ws.Range("B15").Value = 45954.13541666666
$ws.Range("A16").Value = 5328
$ws.Range("B16").Value = 45954.14583333334
$ws.Range("A17").Value = 5357
$ws.Range("B17").Value = 45954.15625
$ws.Range("A18").Value = 5418
$ws.Range("B18").Value = 45954.16666666666
$ws.Range("A19").Value = 5499
$ws.Range("B19").Value = 45954.17708333334
$ws.Range("A20").Value = 5526
$ws.Range("B20").Value = 45954.1875
$ws.Range("A21").Value = 5631
$ws.Range("B21").Value = 45954.19791666666
$ws.Range("A22").Value = 5783
$ws.Range("B22").Value = 45954.20833333334
$ws.Range("A23").Value = 5957
$ws.Range("B23").Value = 45954.21875
$ws.Range("A24").Value = 6136
$ws.Range("B24").Value = 45954.22916666666
$ws.Range("A25").Value = 6263
$ws.Range("B25").Value = 45954.23958333334
$ws.Range("A26").Value = 6576
$ws.Range("B26").Value = 45954.25
$ws.Range("A27").Value = 6709
$ws.Range("B27").Value = 45954.26041666666
$ws.Range("A28").Value = 6816
$ws.Range("B28").Value = 45954.27083333334
$ws.Range("A29").Value = 6902
$ws.Range("B29").Value = 45954.28125
$ws.Range("A30").Value = 7084
$ws.Range("B30").Value = 45954.29166666666
$ws.Range("A31").Value = 7162
$ws.Range("B31").Value = 45954.30208333334
